# Add a new booking entry (Rai Davis / Townhouse) to row 5 of Sheet1.
# Row 4 is intentionally left blank, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A5: Date column -> stored as a real date serial, formatted like a short date (mm-dd-yy -> builtin numFmtId 14)
$ws.Range("A5").Value = 45854
$ws.Range("A5").NumberFormat = "mm-dd-yy"

# B5: Guest name
$ws.Range("B5").Value = "Rai Davis"

# C5 / D5: Check-in / Check-out -> real dates formatted as "d-mmm" (builtin numFmtId 16)
$ws.Range("C5").Value = 45842
$ws.Range("C5").NumberFormat = "d-mmm"

$ws.Range("D5").Value = 45846
$ws.Range("D5").NumberFormat = "d-mmm"

# E5 (Nights) intentionally left blank

# F5: Room/Unit
$ws.Range("F5").Value = "Townhouse"

# G5 / H5: Earnings / Expenses
$ws.Range("G5").Value = 318.48
$ws.Range("H5").Value = 275.99

# I5: Notes
$ws.Range("I5").Value = "Rai Davis pdf"

# J5: Month
$ws.Range("J5").Value = "July"

# Autofit the columns so the new, wider content is fully visible.
$ws.Range("A1:J5").EntireColumn.AutoFit()

# Leave the selection the way it was left in the saved workbook.
$ws.Range("J10").Select() | Out-Null
